$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 847, shifting the existing row 847 (and everything
# below it) down to row 848. This mirrors a new daily price record being
# added to the consolidated sheet.
$ws.Range("A847:T847").EntireRow.Insert()

# Populate the newly inserted row 847 with the new record's data. The
# columns that are constant for every record in this sheet (market,
# region, product taxonomy, unit, kg/unit) are copied straight from the
# neighbouring row.
$ws.Cells.Item(847, 1).Value = $ws.Cells.Item(848, 1).Value()        # Mercado ID
$ws.Cells.Item(847, 2).Value = $ws.Cells.Item(848, 2).Value()        # Mercado
$ws.Cells.Item(847, 3).Value = $ws.Cells.Item(848, 3).Value()        # Region

$dCell = $ws.Cells.Item(847, 4)
$dCell.Value = 45077                                                 # Fecha
$dCell.NumberFormat = $ws.Cells.Item(848, 4).NumberFormat

$ws.Cells.Item(847, 5).Value = $ws.Cells.Item(848, 5).Value()        # Codreg
$ws.Cells.Item(847, 6).Value = $ws.Cells.Item(848, 6).Value()        # Tipo
$ws.Cells.Item(847, 7).Value = $ws.Cells.Item(848, 7).Value()        # Producto ID
$ws.Cells.Item(847, 8).Value = $ws.Cells.Item(848, 8).Value()        # Producto
$ws.Cells.Item(847, 9).Value = $ws.Cells.Item(848, 9).Value()        # Categoria ID
$ws.Cells.Item(847, 10).Value = $ws.Cells.Item(848, 10).Value()      # Categoria

$ws.Cells.Item(847, 11).Value = "Fukumoto"                           # Variedad
$ws.Cells.Item(847, 12).Value = "Primera"                            # Calidad
$ws.Cells.Item(847, 13).Value = 14                                   # Volumen
$ws.Cells.Item(847, 14).Value = 235000                               # Precio minimo
$ws.Cells.Item(847, 15).Value = 240000                               # Precio maximo
$ws.Cells.Item(847, 16).Value = 237500                               # Precio promedio ponderado

$ws.Cells.Item(847, 17).Value = $ws.Cells.Item(848, 17).Value()      # Unidad de comercializacion

$ws.Cells.Item(847, 18).Value = "Provincia de Melipilla"             # Origen
$ws.Cells.Item(847, 19).Value = 594                                  # Precio $/Kg

$ws.Cells.Item(847, 20).Value = $ws.Cells.Item(848, 20).Value()      # Kg / unidad
